# Commit: "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
#
# Semantic changes applied to sheet "Hoja1":
#   1. The "Periodo Mora" value shown in E16, E17 and E18 changes from
#      "2508" to "2509" (all three cells share the same text).
#   2. Those same three cells (E16, E17, E18) pick up a center horizontal
#      alignment that they did not have before.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$cells = @("E16", "E17", "E18")

foreach ($addr in $cells) {
    $cell = $ws.Range($addr)

    # Update the period value: "2508" -> "2509".
    $cell.Value = "2509"

    # Apply center horizontal alignment to the cell.
    $cell.HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter
}
